$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column G ("Recorded By") contains values such as "dnasr281@gmail.com, System".
# Swap the order of the two names to "System, dnasr281@gmail.com" for every
# cell that holds exactly that text (whole-cell match only, so cells that
# contain just "System" or just "dnasr281@gmail.com" are left untouched).
$col = $ws.Columns.Item(7)
$xlWhole = 1
$col.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com", $xlWhole) | Out-Null
